$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1: "I am enclosing ..." -> drop <w:ind w:left="720"/> ---
$p1xml = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/><w:color w:val=`"000000`"/></w:rPr><w:t>I am</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:t xml:space=`"preserve`"> enclosing the following documents towards proof of my identification and address</w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(1).Range.InsertXML($p1xml)

# --- Paragraph 2: "{#summary1}...{title}{/summary1}" -> drop ind, drop leading tab, drop tab-only run ---
$p2xml = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:t>{#summary1}</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:br/></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/><w:b/><w:bCs/></w:rPr><w:t>{title}</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:t>{/summary1}</w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(2).Range.InsertXML($p2xml)

# --- Paragraph 3: "Yours faithfully," -> drop pStyle/tabs/ind/contextualSpacing/jc, strip leading spaces ---
$p3xml = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:t>Yours faithfully,</w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(3).Range.InsertXML($p3xml)

# --- Paragraph 4: "_______________________" -> drop ind, append softHyphen + "____" run ---
$p4xml = "<w:p $wns>" +
  "<w:pPr><w:contextualSpacing/><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:t>_______________________</w:t></w:r>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:softHyphen/><w:t>____</w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(4).Range.InsertXML($p4xml)

# --- Paragraph 5: "{name}" -> drop ind, split into {  + fhnameInPancardExactSpelling (proofErr) + } ---
$p5xml = "<w:p $wns>" +
  "<w:pPr><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr></w:pPr>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:t>{</w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:t>fhnameInPancardExactSpelling</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:rPr><w:rFonts w:ascii=`"Tahoma`" w:hAnsi=`"Tahoma`" w:cs=`"Tahoma`"/></w:rPr><w:t>}</w:t></w:r>" +
  "</w:p>"
$d.Paragraphs(5).Range.InsertXML($p5xml)
